$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terminplan")

$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 2

$ws.Range("G15").Select()
